$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("addCustomerTest")

# Update D6 to the new postal code value
$ws.Range("D6").Value = "A1B2Z1"

# Move the active selection from C5 to A6
$ws.Range("A6").Select()
